$wb = $excel.ActiveWorkbook

# --- Update selections on the existing sheets -------------------------------
$wsTest = $wb.Worksheets.Item("Test")
[void]$wsTest.Range("B15").Select()

$wsConversion = $wb.Worksheets.Item("Conversion")
[void]$wsConversion.Range("E11").Select()

# --- Add "Multiple" worksheet (after "Conversion") ---------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMultiple = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$wsMultiple.Name = "Multiple"
$wsMultiple.PageSetup.TopMargin = 56.692913399999995
$wsMultiple.PageSetup.BottomMargin = 56.692913399999995

# Header strings first (so shared-string table order is A,B,C,D,E,F,...)
$wsMultiple.Range("C6").Value = "A"
$wsMultiple.Range("D6").Value = "B"
$wsMultiple.Range("E6").Value = "C"
$wsMultiple.Range("I10").Value = "D"
$wsMultiple.Range("J10").Value = "E"
$wsMultiple.Range("K10").Value = "F"

# First data block (rows 7-9)
$wsMultiple.Range("C7").Value = 1
$wsMultiple.Range("D7").Value = "a"
$wsMultiple.Range("E7").Value = $true

$wsMultiple.Range("C8").Value = 2
$wsMultiple.Range("D8").Value = "b"
$wsMultiple.Range("E8").Value = $true

$wsMultiple.Range("C9").Value = 3
$wsMultiple.Range("D9").Value = "c"
$wsMultiple.Range("E9").Value = $false

# Second data block (rows 11-13)
$wsMultiple.Range("I11").Value = 4
$wsMultiple.Range("J11").Value = "d"
$wsMultiple.Range("K11").Value = $false

$wsMultiple.Range("I12").Value = 5
$wsMultiple.Range("J12").Value = "e"
$wsMultiple.Range("K12").Value = $true

$wsMultiple.Range("I13").Value = 6
$wsMultiple.Range("J13").Value = "f"
$wsMultiple.Range("K13").Value = $true

[void]$wsMultiple.Range("I14").Select()

# --- Add "VariableNames" worksheet (after "Multiple") ------------------------
$wsVariableNames = $wb.Worksheets.Add([Type]::Missing, $wsMultiple)
$wsVariableNames.Name = "VariableNames"
$wsVariableNames.PageSetup.TopMargin = 56.692913399999995
$wsVariableNames.PageSetup.BottomMargin = 56.692913399999995

$wsVariableNames.Range("D5").Value = "With whitespace"
$wsVariableNames.Range("E5").Value = "And some other funky characters: _=?^~!`$@#%§"

$wsVariableNames.Range("D6").Value = 1
$wsVariableNames.Range("E6").Value = "a"

$wsVariableNames.Range("D7").Value = 2
$wsVariableNames.Range("E7").Value = "b"

$wsVariableNames.Range("D8").Value = 3
$wsVariableNames.Range("E8").Value = "c"

$wsVariableNames.Range("D9").Value = 4
$wsVariableNames.Range("E9").Value = "d"

$wsVariableNames.Columns.Item(4).ColumnWidth = 15.0
$wsVariableNames.Columns.Item(5).ColumnWidth = 43.166666666666664

[void]$wsVariableNames.Range("D12").Select()

# --- Defined names -----------------------------------------------------------
$wb.Names.Add("AAA", "=Multiple!`$C`$6:`$E`$9")
$wb.Names.Add("BBB", "=Multiple!`$I`$10:`$K`$13")
$wb.Names.Add("VariableNames", "=VariableNames!`$D`$5:`$E`$9")
